# edit.ps1 -- apply the 2025-12-23 Betfair odds-sheet update
# Summary of the change:
#  1) Row 2 (Australian A-League Men) and Row 3 (Algerian Ligue 1): several odds refreshed.
#  2) A new fixture "Friendly Matches: Serra Branca EC v Maguary" is inserted as the new row 4,
#     pushing the existing "Portuguese Primeira Liga: Guimaraes v Sporting Lisbon" row down to
#     row 5 (its odds are refreshed too).
#  3) A new fixture "Friendly Matches: Necaxa v Atletico San Luis" is appended as row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a date/time-looking string into a cell as literal TEXT,
# avoiding Excel's automatic "looks like a date" -> date-serial conversion.
# We stage the text in an out-of-range scratch cell formatted as Text, copy
# it, and paste-special (values only) into the destination -- the pasted
# cell keeps the destination's original (General) number format, so no new
# cell style is introduced, but the content remains a plain string.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("BZ1000")

function Set-TextValue($rangeAddress, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)
}

# --- Step 1: refresh odds already present in row 2 (Australian A-League Men) ---
$ws.Range("F2").Value = 1.72
$ws.Range("G2").Value = 1.73
$ws.Range("J2").Value = 4.1
$ws.Range("K2").Value = 4.2
$ws.Range("N2").Value = 4.4
$ws.Range("O2").Value = 1.27
$ws.Range("P2").Value = 2.14
$ws.Range("R2").Value = 1.45
$ws.Range("T2").Value = 1.82
$ws.Range("U2").Value = 2.14
$ws.Range("W2").Value = 2.36
$ws.Range("AA2").Value = 140
$ws.Range("AC2").Value = 9.199999999999999
$ws.Range("AE2").Value = 70
$ws.Range("AH2").Value = 19.5
$ws.Range("AI2").Value = 70
$ws.Range("AM2").Value = 100
$ws.Range("AO2").Value = 70

# --- Step 2: refresh odds already present in row 3 (Algerian Ligue 1) ---
$ws.Range("N3").Value = 2.9
$ws.Range("S3").Value = 4.4
$ws.Range("U3").Value = 1.52

# --- Step 3: insert a blank row at position 4; this pushes the existing
#     "Portuguese Primeira Liga" row (and everything below it) down by one ---
$ws.Rows.Item(4).Insert()

# --- Step 4: populate the new row 4 ("Friendly Matches": Serra Branca EC v Maguary) ---
$ws.Range("A4").Value = "Friendly Matches"
Set-TextValue "B4" "2025-12-23"
$ws.Range("C4").Value = "16:00:00"
$ws.Range("D4").Value = "Serra Branca EC"
$ws.Range("E4").Value = "Maguary"
$ws.Range("F4").Value = 1.04
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 1.04
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 1.02
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 1.25
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 1.24
$ws.Range("Q4").Value = 1.29
$ws.Range("R4").Value = 1.18
$ws.Range("S4").Value = 1.29
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("V4").Value = 1.01
$ws.Range("W4").Value = 1.01
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# --- Step 5: refresh the odds for row 5 (the shifted-down Portuguese Primeira
#     Liga fixture). Columns A-E already hold the right League/Date/Time/Home/
#     Away values because they were carried down by the row insert above, so
#     only the odds columns (F onward) need to be (re)written. ---
$ws.Range("F5").Value = 8.2
$ws.Range("G5").Value = 8.8
$ws.Range("H5").Value = 1.46
$ws.Range("I5").Value = 1.47
$ws.Range("J5").Value = 4.9
$ws.Range("K5").Value = 5.1
$ws.Range("L5").Value = 1.38
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 3.9
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 2.04
$ws.Range("Q5").Value = 1.92
$ws.Range("R5").Value = 1.37
$ws.Range("S5").Value = 3.4
$ws.Range("T5").Value = 2.2
$ws.Range("U5").Value = 1.76
$ws.Range("V5").Value = 3.1
$ws.Range("W5").Value = 1.13
$ws.Range("X5").Value = 16.5
$ws.Range("Y5").Value = 7.6
$ws.Range("Z5").Value = 7.8
$ws.Range("AA5").Value = 12
$ws.Range("AB5").Value = 25
$ws.Range("AC5").Value = 11
$ws.Range("AD5").Value = 9.8
$ws.Range("AE5").Value = 16
$ws.Range("AF5").Value = 80
$ws.Range("AG5").Value = 34
$ws.Range("AH5").Value = 29
$ws.Range("AI5").Value = 40
$ws.Range("AJ5").Value = 340
$ws.Range("AK5").Value = 160
$ws.Range("AL5").Value = 140
$ws.Range("AM5").Value = 200
$ws.Range("AN5").Value = 250
$ws.Range("AO5").Value = 8

# --- Step 6: append new row 6 ("Friendly Matches": Necaxa v Atletico San Luis) ---
$ws.Range("A6").Value = "Friendly Matches"
Set-TextValue "B6" "2025-12-23"
$ws.Range("C6").Value = "18:00:00"
$ws.Range("D6").Value = "Necaxa"
$ws.Range("E6").Value = "Atletico San Luis"
$ws.Range("F6").Value = 1.04
$ws.Range("G6").Value = 1000
$ws.Range("H6").Value = 1.04
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 1.02
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 1.25
$ws.Range("O6").Value = 1.02
$ws.Range("P6").Value = 1.24
$ws.Range("Q6").Value = 1.32
$ws.Range("R6").Value = 1.18
$ws.Range("S6").Value = 1.32
$ws.Range("T6").Value = 1.01
$ws.Range("U6").Value = 1.01
$ws.Range("V6").Value = 1.01
$ws.Range("W6").Value = 1.01
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 1000
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000

# --- Step 7: clean up the scratch cell used for text-staging above ---
$scratch.Clear()
